$wb = $excel.ActiveWorkbook

# Update "展览" sheet (F2, F3, F5 - 想去人数 counts)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10163
$ws1.Range("F3").Value = 228
$ws1.Range("F5").Value = 626

# Update "全部类型" sheet (F2, F3, F5 - 想去人数 counts)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10163
$ws4.Range("F3").Value = 228
$ws4.Range("F5").Value = 626
